$d = $word.ActiveDocument

# 1. Merge the split "La La Land Soundtrack" runs (removes the spell-check
#    proofErr markup) into a single clean run by replacing the text with
#    itself via Find/Replace.
$d.Content.Find.Execute("La La Land Soundtrack", $true, $false, $false, $false, $false, $true, 1, $false, "La La Land Soundtrack", 2) | Out-Null

# 2. Append two more entries to the "Happy Songs" list, after
#    "Hayley Kiyoko Album" (which currently hosts the _GoBack bookmark).
#    InsertParagraphAfter on that paragraph's range creates a new list
#    paragraph inheriting the ListParagraph style / numPr, and carries the
#    bookmark along with the paragraph mark, leaving it as the very last
#    thing in the document - matching the target layout.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1.Range.Text = "Happy – Pharrell Williams"

$lastPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara2.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.Text = "Lady is a Tramp – Ella Fitzgerald"
